
# New weekly price-report row for "Vega Modelo de Temuco" (Maracuyá) needs to be
# inserted right after the current row 54 (i.e. at row 55), pushing every
# following record down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 55; it inherits formatting (the date
# style in column D) from the row above, and every row below shifts down.
$ws.Rows.Item(55).Insert()

# Populate the new row with the latest weekly quote.
$ws.Cells.Item(55, 1).Value = 10
$ws.Cells.Item(55, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(55, 3).Value = "La Araucanía"
$ws.Cells.Item(55, 4).Value = 45174
$ws.Cells.Item(55, 5).Value = 9
$ws.Cells.Item(55, 6).Value = "Fruta"
$ws.Cells.Item(55, 7).Value = 100108
$ws.Cells.Item(55, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(55, 9).Value = 100108003
$ws.Cells.Item(55, 10).Value = "Maracuyá"
$ws.Cells.Item(55, 11).Value = "Sin especificar"
$ws.Cells.Item(55, 12).Value = "Primera"
$ws.Cells.Item(55, 13).Value = 50
$ws.Cells.Item(55, 14).Value = 50000
$ws.Cells.Item(55, 15).Value = 50000
$ws.Cells.Item(55, 16).Value = 50000
$ws.Cells.Item(55, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(55, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(55, 19).Value = 2778
$ws.Cells.Item(55, 20).Value = 18
